$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the "Price" column (D) stores numeric-looking values as TEXT in this
# workbook. Force the Text number format before assignment so the value is
# not silently reinterpreted as a Number by the COM layer.
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Simple price (column D) updates ---
Set-TextValue "D2"  "246.27"
Set-TextValue "D3"  "24.21"
Set-TextValue "D6"  "6.495"
Set-TextValue "D8"  "0.8128"
Set-TextValue "D9"  "0.8606"
Set-TextValue "D11" "0.06942"
Set-TextValue "D12" "0.03147"
Set-TextValue "D15" "3.759"
Set-TextValue "D16" "0.001525"
Set-TextValue "D17" "0.04668"
Set-TextValue "D18" "0.0005974"
Set-TextValue "D19" "0.006137"
Set-TextValue "D20" "0.001237"
Set-TextValue "D21" "0.004631"
Set-TextValue "D24" "2.149"
Set-TextValue "D26" "0.1320"
Set-TextValue "D28" "0.0002333"
Set-TextValue "D40" "0.03690"
Set-TextValue "D45" "0.00005266"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.4402"
Set-TextValue "D48" "0.002381"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"

# --- Row 18 (One / ONE): Volume(1h) label gains a "Worstin24h" suffix ---
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Rows 41-43: the three coins rotate (KickToken, BKEXToken, CEJI each
# shift down one slot), each keeping its own row's rank prefix/suffix text
# but with new price + name + link values ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006253"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1055"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003001"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

# --- Row 44 (LocalTraders / LCT): price update + "Bestin24h" suffix removed ---
Set-TextValue "D44" "0.008502"
$ws.Range("E44").Value = "43LocalTradersLCT"
